$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 438; this shifts existing rows 438:563
# down to 439:564 and extends the used range to A1:R564.
$ws.Rows("438:438").Insert()

# Populate the newly inserted row 438 with the new data record.
$ws.Range("A438").Value = 6
$ws.Range("B438").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C438").Value = "Metropolitana"
$ws.Range("D438").Value = 44736
$ws.Range("E438").Value = 13
$ws.Range("F438").Value = 100112012
$ws.Range("G438").Value = "Espinaca"
$ws.Range("H438").Value = "Sin especificar"
$ws.Range("I438").Value = "Primera"
$ws.Range("J438").Value = 310
$ws.Range("K438").Value = 9000
$ws.Range("L438").Value = 10000
$ws.Range("M438").Value = 9419
$ws.Range("N438").Value = "$/cuna 10 kilos"
$ws.Range("O438").Value = "Región Metropolitana"
$ws.Range("P438").Value = 942
$ws.Range("Q438").Value = 10
$ws.Range("R438").Value = "Hortaliza"
